# 11.6.1.1a — add a new "2023" column (column N) to the waste-disposal
# table, mirroring the existing 2022 column (M), and bump L7's number
# style to match the rest of its row family.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row heights: header block gets a bit taller / becomes an explicit
# (custom) height instead of the auto-fit default.
$ws.Rows.Item(1).RowHeight = 66
$ws.Rows.Item(2).RowHeight = 14.25
$ws.Rows.Item(3).RowHeight = 14.25

# --- New column N ("2023"), formatted the same as column M ---

# N3: blank divider cell under the thick bottom border, same style as M3.
$ws.Range("M3").Copy()
$ws.Range("N3").PasteSpecial(-4122)

# N4: year label.
$ws.Range("M4").Copy()
$ws.Range("N4").PasteSpecial(-4122)
$ws.Range("N4").Value = 2023

# N5: "Garbage chute" 2023 value.
$ws.Range("M5").Copy()
$ws.Range("N5").PasteSpecial(-4122)
$ws.Range("N5").Value = 0

# N6: "Collection by truck, container" 2023 value.
$ws.Range("M6").Copy()
$ws.Range("N6").PasteSpecial(-4122)
$ws.Range("N6").Value = 48.5

# N7: "Dumping into garbage heaps" 2023 value.
$ws.Range("M7").Copy()
$ws.Range("N7").PasteSpecial(-4122)
$ws.Range("N7").Value = 23.2

# N8: "Burning" 2023 value.
$ws.Range("M8").Copy()
$ws.Range("N8").PasteSpecial(-4122)
$ws.Range("N8").Value = 19.3

# N9: "Instillation" 2023 value.
$ws.Range("M9").Copy()
$ws.Range("N9").PasteSpecial(-4122)
$ws.Range("N9").Value = 9.1

# L7 (2021, "Dumping into garbage heaps") switches to the same number
# style used by the rest of that row (e.g. L8), value is unchanged.
$oldL7 = $ws.Range("L7").Value()
$ws.Range("L8").Copy()
$ws.Range("L7").PasteSpecial(-4122)
$ws.Range("L7").Value = $oldL7

# Reset the active cell away from the old N7 selection.
$ws.Range("A1").Select()
